# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (it already has the right column layout:
#    基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名),
#    place the copy right before "总计", and rename it to "2022-Q1".
# 2. Overwrite the two data rows of the new sheet with the 2022-Q1 holdings.
# 3. Add a new "2022-Q1" row at the top of the "总计" (summary) sheet and
#    push the existing "2021-Q4" / "2020-Q4" rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by copying "2021-Q4"
# ---------------------------------------------------------------------
$src   = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Helper: write a value as a *text* cell (matches the "numeric-looking
# string" cells used throughout this workbook) without leaving behind any
# stray style index.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: 001914 / 中信建投聚利混合A
Set-TextCell $newSheet.Range("B2") "001914"
Set-TextCell $newSheet.Range("C2") "中信建投聚利混合A"
Set-TextCell $newSheet.Range("D2") "0.13"
Set-TextCell $newSheet.Range("E2") "39.07"
Set-TextCell $newSheet.Range("F2") "2.08"
Set-TextCell $newSheet.Range("G2") "0.0027"
$newSheet.Range("H2").Value = 6

# Row 3: 000041 / 华夏全球精选股票(QDII)
Set-TextCell $newSheet.Range("B3") "000041"
Set-TextCell $newSheet.Range("C3") "华夏全球精选股票(QDII)"
Set-TextCell $newSheet.Range("D3") "0.02"
Set-TextCell $newSheet.Range("E3") "39.07"
Set-TextCell $newSheet.Range("F3") "2.08"
Set-TextCell $newSheet.Range("G3") "0.0004"
$newSheet.Range("H3").Value = 6

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet - insert a 2022-Q1 row on top,
# push old rows down (row2 -> row3, row3 -> row4).
# ---------------------------------------------------------------------
$sumWs = $wb.Worksheets.Item("总计")

# Push 2020-Q4 row (was row 2) down to row 4
$sumWs.Range("A4").Value = 2
$sumWs.Range("A3").Copy()
$sumWs.Range("A4").PasteSpecial(-4122)   # xlPasteFormats - keep the "A" column style
$sumWs.Range("B4").Value = "2020-Q4"
$sumWs.Range("C4").Value = 2
$sumWs.Range("D4").Value = 0.01

# Push 2021-Q4 row (was row 2) down to row 3
$sumWs.Range("A3").Value = 1
$sumWs.Range("B3").Value = "2021-Q4"
$sumWs.Range("C3").Value = 2
$sumWs.Range("D3").Value = 0.13

# New 2022-Q1 row at row 2
$sumWs.Range("A2").Value = 0
$sumWs.Range("B2").Value = "2022-Q1"
$sumWs.Range("C2").Value = 2
$sumWs.Range("D2").Value = 0
